# Apply the edits described by the diff:
#  - remove the extra "Примечания.1" / "Финальное решение.1" columns (P:Q)
#  - remove the now-unused sample data rows (old rows 5:8)
#  - clear the remaining sample rows (2:4) back to "empty" template rows,
#    keeping the styled A/B columns
#  - widen columns N and O (new <cols> entries in the target sheet)
#  - leave the selection on G9, matching the saved workbook view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing "Примечания.1" / "Финальное решение.1" columns entirely.
$ws.Columns("P:Q").Delete() | Out-Null

# Drop the extra sample rows, leaving just the 3 template rows (2:4).
$ws.Rows("5:8").Delete() | Out-Null

# Clear out the remaining sample data (values only; keep the A/B styles).
$ws.Range("A2:O4").ClearContents() | Out-Null

# New explicit column widths for N (15.90625) and O (20 characters).
$ws.Columns("N").ColumnWidth = 15
$ws.Columns("O").ColumnWidth = 19.1666666

# Restore the saved selection state (cell G9).
$ws.Range("G9").Select() | Out-Null
